$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BNB)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.12%'

# Row 3 (OKB)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '36.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-4.50%'

# Row 4 (HuobiToken)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.116'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.50%'

# Row 5 (Cronos)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07717'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-2.48%'

# Row 6 (KuCoinToken)
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.395'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.57%'

# Row 7 (FTXToken)
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.292'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.29%'

# Row 8 (BTSEToken)
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.844'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.99%'

# Row 9 (MXToken)
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.965'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.36%'

# Row 10 (LiechtensteinCryptoassetsExchange)
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9187'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.23%'

# Row 11 (WazirX)
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1094'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-10.76%'

# Row 12 (MandalaExchangeToken)
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1854'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.10%'

# Row 13 (BitrueCoin)
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08773'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.88%'

# Row 14 (BitMartToken)
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03332'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.63%'

# Row 15 (BitForexToken)
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09512'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.26%'

# Row 16 (TigerCash)
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001377'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.10%'

# Row 17 (LEO)
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006115'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '6.25%'

# Row 18 (GateToken)
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.361'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-4.77%'

# Row 19 (BitpandaEcosystemToken)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3448'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.15%'

# Row 20 (MCDex)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.303'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '19.74%'

# Row 21 (ProBitToken)
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.47%'

# Row 22 (ZBToken)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2314'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-10.74%'

# Row 23 (CoinExToken)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04321'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.97%'

# Row 24 (BitKan)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001199'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-4.01%'

# Row 25 (HotbitToken)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004250'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.24%'

# Row 26 (NitroEx)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001326'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '8.70%'

# Row 27 (UpBots)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002903'

# Row 39 (One)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02069'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.77%'

# Row 40 (IDEX)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04934'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-4.62%'

# Row 41 (KickToken)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007485'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.06%'

# Row 42 (BKEXToken)
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.88%'

# Row 43 (Dexo)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008544'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.18%'

# Row 44 (CEJI)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002063'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2.68%'

# Row 45 (LocalTraders)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008348'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-3.15%'

# Row 46 (CoinLion)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006324'

# Row 47 (Kangarootoken)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.15%'

# Row 48 (BOLO)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002841'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-14.13%'

# Row 49 (CoinbaseStockToken)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001442'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '20.03%'

# Row 50 (CryptobidCoin)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.15%'

# Row 51 (SpecialPowerGold)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.15%'
